$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 343, pushing existing rows 343:370 down to 344:371
$ws.Rows("343:343").Insert()

# Populate the new row 343 with this week's data (new price record)
$ws.Range("A343").Value = 7
$ws.Range("B343").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C343").Value = "Ñuble"
$ws.Range("D343").Value = 45166
$ws.Range("E343").Value = 16
$ws.Range("F343").Value = 100112043
$ws.Range("G343").Value = "Pepino ensalada"
$ws.Range("H343").Value = "Sin especificar"
$ws.Range("I343").Value = "Primera"
$ws.Range("J343").Value = 100
$ws.Range("K343").Value = 12000
$ws.Range("L343").Value = 12000
$ws.Range("M343").Value = 12000
$ws.Range("N343").Value = "`$/caja 60 unidades"
$ws.Range("O343").Value = "Región de Arica y Parinacota"
$ws.Range("P343").Value = 200
$ws.Range("Q343").Value = 60
$ws.Range("R343").Value = "Hortaliza"
